$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full replacement data block for rows 2-9 (A:AQ), reflecting the
# refreshed capital-structure figures and the reordered/updated company list.
$arr = New-Object 'object[,]' 8,43

$arr[0,0] = 'Israel'
$arr[0,1] = '''7'
$arr[0,2] = 'Bank (Money Center)'
$arr[0,3] = 0.013
$arr[0,4] = 0.0708
$arr[0,5] = $null
$arr[0,6] = 0
$arr[0,7] = 0
$arr[0,8] = -0.009576479694674596
$arr[0,9] = -0.006064048255648862
$arr[0,10] = 1825
$arr[0,11] = 0.1448389708099871
$arr[0,12] = 839.4300000000001
$arr[0,13] = 0.02621138221536656
$arr[0,14] = 0.4599616438356165
$arr[0,15] = 707.03
$arr[0,16] = 0.02207716375127243
$arr[0,17] = 0.387413698630137
$arr[0,18] = 132.4
$arr[0,19] = 0.1577260760277808
$arr[0,20] = 134568.6
$arr[0,21] = 4.201933465311909
$arr[0,22] = 0.08033033033033034
$arr[0,23] = 0.05973061924931371
$arr[0,24] = 0.02059971108101663
$arr[0,25] = -9.035820129722952
$arr[0,26] = 0
$arr[0,27] = 0.03817823786333101
$arr[0,28] = -0.03841709948549001
$arr[0,29] = 39791.2
$arr[0,30] = 1435.827797244194
$arr[0,31] = 41227.02779724419
$arr[0,32] = -93341.57220275581
$arr[0,33] = 0.5628076643596949
$arr[0,34] = 0.5109505210167671
$arr[0,35] = 1.522299400786154
$arr[0,36] = 1.732343640923264
$arr[0,37] = 0
$arr[0,38] = 0
$arr[0,39] = 238.9861861861862
$arr[0,40] = $null
$arr[0,41] = -560.6100432597947
$arr[0,42] = $null
$arr[1,0] = 'Israel'
$arr[1,1] = 'First International Bank of Israel Ltd (TASE:FIBI)'
$arr[1,2] = 'Bank (Money Center)'
$arr[1,3] = 0.013
$arr[1,4] = 0.149
$arr[1,5] = $null
$arr[1,6] = 0
$arr[1,7] = 0
$arr[1,8] = -0.004762926778355925
$arr[1,9] = -0.003242341940290648
$arr[1,10] = 222.7
$arr[1,11] = 0.206720504966119
$arr[1,12] = 68.7
$arr[1,13] = 0.02585525572993113
$arr[1,14] = 0.308486753480018
$arr[1,15] = 68.7
$arr[1,16] = 0.02585525572993113
$arr[1,17] = 0.308486753480018
$arr[1,18] = 0
$arr[1,19] = 0
$arr[1,20] = 15192.7
$arr[1,21] = 5.717775017876633
$arr[1,22] = 0.09157072368421053
$arr[1,23] = 0.04944323687223476
$arr[1,24] = 0.04212748681197577
$arr[1,25] = -0.1841209720653249
$arr[1,26] = 0.0005969831498144858
$arr[1,27] = 0.03698964490644516
$arr[1,28] = -0.03639266175663068
$arr[1,29] = 1768.7
$arr[1,30] = 131.1555050916142
$arr[1,31] = 1899.855505091614
$arr[1,32] = -13292.84449490839
$arr[1,33] = 0.4169133323704505
$arr[1,34] = 0.4107302361846776
$arr[1,35] = 1.24982736293374
$arr[1,36] = 1.257941017207945
$arr[1,37] = 0
$arr[1,38] = 0
$arr[1,39] = 83.82464454976304
$arr[1,40] = $null
$arr[1,41] = -629.9926300904448
$arr[1,42] = $null
$arr[2,0] = 'Israel'
$arr[2,1] = 'F.I.B.I. Holdings Ltd (TASE:FIBIH)'
$arr[2,2] = 'Bank (Money Center)'
$arr[2,3] = 0.0129
$arr[2,4] = 0.15
$arr[2,5] = $null
$arr[2,6] = 0
$arr[2,7] = 0
$arr[2,8] = 0
$arr[2,9] = 0
$arr[2,10] = 106.4
$arr[2,11] = 0.09876543209876544
$arr[2,12] = 31.6
$arr[2,13] = 0.02941997951773578
$arr[2,14] = 0.2969924812030075
$arr[2,15] = 31.6
$arr[2,16] = 0.02941997951773578
$arr[2,17] = 0.2969924812030075
$arr[2,18] = 0
$arr[2,19] = 0
$arr[2,20] = 15195
$arr[2,21] = 14.14672749278466
$arr[2,22] = 0.09039164047234731
$arr[2,23] = 0.06801602511410608
$arr[2,24] = 0.02237561535824123
$arr[2,25] = -0.1487983425414365
$arr[2,26] = 0
$arr[2,27] = 0.03787587555202758
$arr[2,28] = -0.03787587555202758
$arr[2,29] = 1768.7
$arr[2,30] = 0
$arr[2,31] = 1768.7
$arr[2,32] = -13426.3
$arr[2,33] = 0.622168284789644
$arr[2,34] = 0.3932541800071149
$arr[2,35] = 1.086956169751137
$arr[2,36] = 1.255099369940359
$arr[2,37] = 0
$arr[2,38] = 0
$arr[2,39] = $null
$arr[2,40] = $null
$arr[2,41] = $null
$arr[2,42] = $null
$arr[3,0] = 'Israel'
$arr[3,1] = 'Bank Leumi le- Israel B.M. (TASE:LUMI)'
$arr[3,2] = 'Bank (Money Center)'
$arr[3,3] = -0.0361
$arr[3,4] = -0.0371
$arr[3,5] = $null
$arr[3,6] = 0
$arr[3,7] = 0
$arr[3,8] = 0
$arr[3,9] = 0
$arr[3,10] = 571
$arr[3,11] = 0.1877363143185928
$arr[3,12] = 308.6
$arr[3,13] = 0.03580129469361239
$arr[3,14] = 0.5404553415061296
$arr[3,15] = 176.2
$arr[3,16] = 0.02044130954314485
$arr[3,17] = 0.3085814360770578
$arr[3,18] = 132.4
$arr[3,19] = 0.4290343486714194
$arr[3,20] = 33036.8
$arr[3,21] = 3.832664330958956
$arr[3,22] = 0.05652457977786138
$arr[3,23] = 0.05655880939136944
$arr[3,24] = -0.00003422961350806053
$arr[3,25] = 1.017428246470864
$arr[3,26] = 0
$arr[3,27] = 0.03817823786333101
$arr[3,28] = -0.03817823786333101
$arr[3,29] = 9240
$arr[3,30] = 0
$arr[3,31] = 9240
$arr[3,32] = -23796.8
$arr[3,33] = 0.5173630163831622
$arr[3,34] = 0.4569326172744266
$arr[3,35] = 1.567951505567635
$arr[3,36] = 1.856948888021849
$arr[3,37] = 0
$arr[3,38] = 0
$arr[3,39] = $null
$arr[3,40] = $null
$arr[3,41] = $null
$arr[3,42] = $null
$arr[4,0] = 'Israel'
$arr[4,1] = 'Bank of Jerusalem Ltd. (TASE:JBNK)'
$arr[4,2] = 'Bank (Money Center)'
$arr[4,3] = 0.034
$arr[4,4] = 0.0519
$arr[4,5] = $null
$arr[4,6] = 0
$arr[4,7] = 0
$arr[4,8] = 0
$arr[4,9] = 0
$arr[4,10] = 21.4
$arr[4,11] = 0.1394136807817589
$arr[4,12] = 4.53
$arr[4,13] = 0.02658450704225352
$arr[4,14] = 0.2116822429906542
$arr[4,15] = 4.53
$arr[4,16] = 0.02658450704225352
$arr[4,17] = 0.2116822429906542
$arr[4,18] = 0
$arr[4,19] = 0
$arr[4,20] = 822.1
$arr[4,21] = 4.824530516431925
$arr[4,22] = 0.08033033033033034
$arr[4,23] = 0.09421053552894694
$arr[4,24] = -0.0138802051986166
$arr[4,25] = 2.260677466863035
$arr[4,26] = 0
$arr[4,27] = 0.03841709948549001
$arr[4,28] = -0.03841709948549001
$arr[4,29] = 504.5
$arr[4,30] = 0
$arr[4,31] = 504.5
$arr[4,32] = -317.6
$arr[4,33] = 0.7475181508371611
$arr[4,34] = 0.6367537548908242
$arr[4,35] = 2.157608695652174
$arr[4,36] = 10.65771812080537
$arr[4,37] = 0
$arr[4,38] = 0
$arr[4,39] = $null
$arr[4,40] = $null
$arr[4,41] = $null
$arr[4,42] = $null
$arr[5,0] = 'Israel'
$arr[5,1] = 'Bank Hapoalim B.M. (TASE:POLI)'
$arr[5,2] = 'Bank (Money Center)'
$arr[5,3] = -0.0579
$arr[5,4] = -0.297
$arr[5,5] = $null
$arr[5,6] = 0
$arr[5,7] = 0
$arr[5,8] = -0.009504434685820186
$arr[5,9] = -0.004752217342910093
$arr[5,10] = 149.6
$arr[5,11] = 0.04999832893285652
$arr[5,12] = 292.2
$arr[5,13] = 0.03197671237374014
$arr[5,14] = 1.953208556149733
$arr[5,15] = 292.2
$arr[5,16] = 0.03197671237374014
$arr[5,17] = 1.953208556149733
$arr[5,18] = 0
$arr[5,19] = 0
$arr[5,20] = 37152.5
$arr[5,21] = 4.065759091257291
$arr[5,22] = 0.01327088212334114
$arr[5,23] = 0.05973061924931371
$arr[5,24] = -0.04645973712597257
$arr[5,25] = 0.8719024868693205
$arr[5,26] = -0.004143470119426824
$arr[5,27] = 0.03837707113295379
$arr[5,28] = -0.04252054125238062
$arr[5,29] = 10894.1
$arr[5,30] = 355.1910951172129
$arr[5,31] = 11249.29109511721
$arr[5,32] = -25903.20890488279
$arr[5,33] = 0.551782295198648
$arr[5,34] = 0.496776075481569
$arr[5,35] = 1.545048114045703
$arr[5,36] = 1.785454339058111
$arr[5,37] = 0
$arr[5,38] = 0
$arr[5,39] = 255.7300469483568
$arr[5,40] = $null
$arr[5,41] = -608.0565470629762
$arr[5,42] = $null
$arr[6,0] = 'Israel'
$arr[6,1] = 'Israel Discount Bank Limited (TASE:DSCT)'
$arr[6,2] = 'Bank (Money Center)'
$arr[6,3] = 0.0257
$arr[6,4] = 0.0834
$arr[6,5] = $null
$arr[6,6] = 0
$arr[6,7] = 0
$arr[6,8] = -0.01001262110394765
$arr[6,9] = -0.006592551447737535
$arr[6,10] = 302.7
$arr[6,11] = 0.1272490331259458
$arr[6,12] = 33
$arr[6,13] = 0.007373807342524523
$arr[6,14] = 0.1090188305252725
$arr[6,15] = 33
$arr[6,16] = 0.007373807342524523
$arr[6,17] = 0.1090188305252725
$arr[6,18] = 0
$arr[6,19] = 0
$arr[6,20] = 10896.8
$arr[6,21] = 2.434875874243067
$arr[6,22] = 0.05705615139577403
$arr[6,23] = 0.05866571128541639
$arr[6,24] = -0.001609559889642363
$arr[6,25] = 0.6019093518956878
$arr[6,26] = -0.003968118369246679
$arr[6,27] = 0.03964816701329115
$arr[6,28] = -0.04361628538253783
$arr[6,29] = 4938.2
$arr[6,30] = 332.0901154103533
$arr[6,31] = 5270.290115410353
$arr[6,32] = -5626.509884589646
$arr[6,33] = 0.5407871717359252
$arr[6,34] = 0.477970583323064
$arr[6,35] = 4.887475307419932
$arr[6,36] = -43.41773959204377
$arr[6,37] = 0
$arr[6,38] = 0
$arr[6,39] = 115.9201877934272
$arr[6,40] = $null
$arr[6,41] = -132.0776968213532
$arr[6,42] = $null
$arr[7,0] = 'Israel'
$arr[7,1] = 'Mizrahi Tefahot Bank Ltd. (TASE:MZTF)'
$arr[7,2] = 'Bank (Money Center)'
$arr[7,3] = 0.05
$arr[7,4] = 0.0708
$arr[7,5] = $null
$arr[7,6] = 0
$arr[7,7] = 0
$arr[7,8] = -0.03366399761930242
$arr[7,9] = -0.02199491960180667
$arr[7,10] = 451.2
$arr[7,11] = 0.2400383039846784
$arr[7,12] = 100.8
$arr[7,13] = 0.01711142798940721
$arr[7,14] = 0.2234042553191489
$arr[7,15] = 100.8
$arr[7,16] = 0.01711142798940721
$arr[7,17] = 0.2234042553191489
$arr[7,18] = 0
$arr[7,19] = 0
$arr[7,20] = 22272.7
$arr[7,21] = 3.780929585115774
$arr[7,22] = 0.09963344079848076
$arr[7,23] = 0.07341066587972275
$arr[7,24] = 0.026222774918758
$arr[7,25] = 1.497183076415849
$arr[7,26] = -0.03293042139495216
$arr[7,27] = 0.03802720344594548
$arr[7,28] = -0.07095762484089764
$arr[7,29] = 10677
$arr[7,30] = 617.3910816250138
$arr[7,31] = 11294.39108162501
$arr[7,32] = -10978.30891837499
$arr[7,33] = 0.6572164969234098
$arr[7,34] = 0.6691508853977813
$arr[7,35] = 2.157894776110111
$arr[7,36] = 2.035278228958202
$arr[7,37] = 0
$arr[7,38] = 0
$arr[7,39] = 177.358803986711
$arr[7,40] = $null
$arr[7,41] = -182.3639355211792
$arr[7,42] = $null

$ws.Range("A2:AQ9").Value = $arr

# Remove the old trailing 10th data row (Bank Hapoalim has moved up into row 7
# and the dataset now only spans through row 9).
$ws.Range("A10:AQ10").Clear()
